$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text changes: row 39 becomes "Unassigned" triple; row 40 becomes "Urophycis sp" triple
$ws.Range("A39").Value = "Unassigned"
$ws.Range("B39").Value = "Unassigned"
$ws.Range("C39").Value = "Unassigned"

$ws.Range("A40").Value = "Urophycis sp"
$ws.Range("B40").Value = "Red White or Spotted hake"
$ws.Range("C40").Value = "Teleost Fish"

# Numeric value changes
$ws.Range("E2").Value = 0.08015624140397205
$ws.Range("F3").Value = 0.4732762888056007
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0.02574075383636235
$ws.Range("D7").Value = 0.02901290799401867
$ws.Range("E8").Value = 0.001283673506812639
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("E10").Value = 0.003062478223395866
$ws.Range("D11").Value = 0.04035681751774634
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("D13").Value = 0.005465701861432426
$ws.Range("E14").Value = 0.001760466523628762
$ws.Range("F14").Value = 0.0411710628668411
$ws.Range("E15").Value = 0.001925510260218958
$ws.Range("E16").Value = 0.0120665309640388
$ws.Range("E17").Value = 0.0458821587720746
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("D19").Value = 0.09575290902528316
$ws.Range("E19").Value = 0.003704314976802186
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("E21").Value = 0.0008618950688599145
$ws.Range("D22").Value = 0.01014076760454444
$ws.Range("E22").Value = 0.1729474977535713
$ws.Range("F22").Value = 0.2076939396082314
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("D25").Value = 0.6262353689348756
$ws.Range("E25").Value = 0.01452384881993728
$ws.Range("E26").Value = 0.0003117492802259266
$ws.Range("D27").Value = 0.001598459978343446
$ws.Range("E27").Value = 0.02416973831398654
$ws.Range("D28").Value = 0.03330984341967309
$ws.Range("E28").Value = 0.01454218701289175
$ws.Range("F28").Value = 0.0004808712255144615
$ws.Range("D29").Value = 0.002165655454529829
$ws.Range("E29").Value = 0.03731822266233886
$ws.Range("F29").Value = 0.0000424298140159819
$ws.Range("E30").Value = 0.006583411270653389
$ws.Range("D31").Value = 0.1302143311390316
$ws.Range("E31").Value = 0.02070381984559241
$ws.Range("E32").Value = 0.003172507381122664
$ws.Range("D33").Value = 0.001289080627696327
$ws.Range("D34").Value = 0.00006875096681047077
$ws.Range("E34").Value = 0.01709119583356256
$ws.Range("D35").Value = 0.002715663189013596
$ws.Range("E35").Value = 0.491188498285379
$ws.Range("F35").Value = 0.2480588360087688
$ws.Range("E36").Value = 0.02099723093286388
$ws.Range("E37").Value = 0.001485393629311768
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0.005464781500430948
$ws.Range("F39").Value = 0.003380241849939891
$ws.Range("D40").Value = 0.02167374228700091
$ws.Range("E40").Value = 0.01879664777832792
$ws.Range("F40").Value = 0.000155575984725267
